$d = $word.ActiveDocument

# Locate the "Вывод" (Conclusion) section's introductory paragraph, a
# duplicate of the "Цель работы" wording ("First Paragraph" style),
# immediately followed by the actual conclusion text ("Body Text" style).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($p.Style.NameLocal -eq "First Paragraph") -and ($t -like "*Получение навыков работы с планировщиками событий cron и at.*")) {
        if ($i -lt $d.Paragraphs.Count) {
            $nxt = $d.Paragraphs.Item($i + 1)
            if ($nxt.Style.NameLocal -eq "Body Text") {
                $targetIndex = $i
            }
        }
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    # Remove the redundant "Получение навыков..." paragraph entirely
    # (including its paragraph mark).
    $target.Range.Delete()
    # The paragraph that used to follow it ("В результате выполнения...")
    # now sits at the same index; promote it to the "First Paragraph"
    # style that the removed paragraph used to carry.
    $promoted = $d.Paragraphs.Item($targetIndex)
    $promoted.Style = "First Paragraph"
}
